$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.154.95"

# Row 3
$ws.Range("D3").Value = "2.271.05"
$ws.Range("E3").Value = "  -1.42%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").Value = "'299.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.32%  "

# Row 6
$ws.Range("D6").Value = "'95.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.36%  "

# Row 7
$ws.Range("E7").Value = "  -2.54%  "

# Row 8
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
$ws.Range("E9").Value = "  -2.62%  "

# Row 10
$ws.Range("D10").Value = "'33.07"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.43%  "

# Row 11
$ws.Range("E11").Value = "  -0.58%  "

# Row 12
$ws.Range("E12").Value = "  -6.79%  "

# Row 13
$ws.Range("E13").Value = "  +0.99%  "

# Row 14
$ws.Range("D14").Value = "'15.85"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.89%  "

# Row 15
$ws.Range("D15").Value = "'6.67"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.09%  "

# Row 16
$ws.Range("D16").Value = "2.623.33"
$ws.Range("E16").Value = "  -1.46%  "

# Row 17
$ws.Range("D17").Value = "2.270.10"
$ws.Range("E17").Value = "  -2.01%  "

# Row 18
$ws.Range("E18").Value = "  -2.37%  "

# Row 19
$ws.Range("D19").Value = "42.126.81"
$ws.Range("E19").Value = "  -1.09%  "

# Row 20
$ws.Range("D20").Value = "'11.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.11%  "

# Row 21
$ws.Range("E21").Value = "  -1.74%  "

# Row 22
$ws.Range("D22").Value = "'5.98"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.38%  "

# Row 23
$ws.Range("D23").Value = "'66.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.65%  "

# Row 24
$ws.Range("D24").Value = "'235.02"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.05%  "

# Row 25
$ws.Range("E25").Value = "  -0.56%  "

# Row 26
$ws.Range("E26").Value = "  +0.12%  "

# Row 27
$ws.Range("D27").Value = "'2.46"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.19%  "

# Row 28
$ws.Range("D28").Value = "'23.73"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.09%  "

# Row 29
$ws.Range("D29").Value = "'2.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.23%  "

# Row 30
$ws.Range("D30").Value = "'167.65"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.25%  "

# Row 31
$ws.Range("D31").Value = "'33.75"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.61%  "

# Row 32
$ws.Range("D32").Value = "'9.15"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.31%  "

# Row 33
$ws.Range("D33").Value = "'1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.01%  "

# Row 34
$ws.Range("D34").Value = "'4.71"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.00%  "

# Row 35
$ws.Range("E35").Value = "  -2.33%  "

# Row 36
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").Value = "'2.37"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.01%  "

# Row 37
$ws.Range("B37").Value = "Celestia"
$ws.Range("C37").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D37").Value = "'16.74"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.13%  "

# Row 38
$ws.Range("D38").Value = "'0.0689"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.09%  "

# Row 39
$ws.Range("D39").Value = "'2.79"
$ws.Range("D39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'0.0985"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.65%  "

# Row 41
$ws.Range("E41").Value = "  -4.49%  "

# Row 42
$ws.Range("E42").Value = "  -2.59%  "

# Row 43
$ws.Range("E43").Value = "  -7.73%  "

# Row 44
$ws.Range("D44").Value = "1.957.55"
$ws.Range("E44").Value = "  -0.54%  "

# Row 45
$ws.Range("E45").Value = "  -1.21%  "

# Row 46
$ws.Range("D46").Value = "'17.58"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.62%  "

# Row 47
$ws.Range("D47").Value = "'9.59"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.19%  "

# Row 48
$ws.Range("E48").Value = "  -4.06%  "

# Row 49
$ws.Range("D49").Value = "2.495.71"
$ws.Range("E49").Value = "  -1.26%  "

# Row 50
$ws.Range("D50").Value = "'52.07"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.78%  "

# Row 51
$ws.Range("E51").Value = "  -2.96%  "
